$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SFR")
$ws.Range("A3:G3").Copy()
$ws.Range("A4:G4").PasteSpecial()
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "SFR3"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "SFR3"
$ws.Range("E4").Value = 0.03
$ws.Range("F4").Value = 0.03
$ws.Range("G4").Value = 3
